# "Colocando header nos graficos"
# Adds a header label to column A (row 1) on each data sheet, matching the
# style already used by the other header cells (B1:E1), removes the
# (accidental) header style from the row-label cells below it, fixes a
# batch of missing Portuguese accents in those row labels, drops the
# now-unused "Teto" row from the emissions sheet, and turns the cost sheet
# into a proper two-row/two-column table with a year header and updated
# figures.

$wb = $excel.ActiveWorkbook

# --- Sheets 1-4: "Fonte/Tecnologia" tables -------------------------------
# Same A-column structure (and the same typos) on every one of these four
# sheets, only the B:E figures differ - so loop over them.
$fonteSheets = 1, 2, 3, 4

foreach ($idx in $fonteSheets) {
    $ws = $wb.Worksheets.Item($idx)

    # New column header, picking up the bold/centered/bordered style that's
    # already applied to B1:E1.
    $ws.Range("A1").Value = "Fonte/Tecnologia"
    $ws.Range("B1").Copy()
    $ws.Range("A1").PasteSpecial(-4122)

    # Row labels A2:A12 should not carry that header style - strip it back
    # off them (text content stays, only the formatting goes).
    $ws.Range("A2:A12").ClearFormats()

    # Spelling/accent fixes.
    $ws.Range("A3").Value = "Gás Natural"
    $ws.Range("A4").Value = "Carvão"
    $ws.Range("A6").Value = "Óleos Comb"
    $ws.Range("A8").Value = "Eólica"
    $ws.Range("A11").Value = "Pot. Compl."
}

# --- Sheet 5: "Emissoes Totais" table -------------------------------------
$ws = $wb.Worksheets.Item(5)

$ws.Range("A1").Value = "Período"
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

$ws.Range("A2:A3").ClearFormats()
$ws.Range("A2").Value = "P.Médio"
$ws.Range("A3").Value = "P.Crítico"

# The "Teto" row is no longer used - remove it entirely (shifts the used
# range back down to A1:E3).
$ws.Range("A4:E4").Delete(-4162)

# --- Sheet 6: "Custo Total" table -----------------------------------------
$ws = $wb.Worksheets.Item(6)

$ws.Range("A1").Value = "Tipo Expansão"
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

# B1 becomes a "2015" year header like the other sheets (keep its existing
# bold/centered/bordered style - only the text changes - so pull the literal
# text value in from a cell that already holds it rather than typing a
# number that would get auto-coerced to a numeric cell).
$wb.Worksheets.Item(1).Range("B1").Copy()
$ws.Range("B1").PasteSpecial(-4163)

$ws.Range("A2:A3").ClearFormats()
$ws.Range("A2").Value = "Expansão Centralizada"
$ws.Range("A3").Value = "Expansão por GD"

$ws.Range("B2").Value = 588
$ws.Range("B3").Value = 99
